$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unstyled) cell used to restore style index after
# forcing NumberFormat = "@" (Text) on D-column price cells below. This keeps
# ambiguous numeric-looking strings (e.g. "1.015") stored as text, matching the
# original inlineStr cell type, without leaving a stray style index behind.
$defaultStyle = $ws.Range("C2").Style

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.094.05"
$ws.Range("D2").Style = $defaultStyle
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.849.04"
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "309.39"
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = "  +0.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").Style = $defaultStyle
$ws.Range("E7").Value = "  +2.71%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3684"
$ws.Range("D8").Style = $defaultStyle
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07245"
$ws.Range("D9").Style = $defaultStyle
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9326"
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.90"
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = "  +2.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07792"
$ws.Range("D12").Style = $defaultStyle
$ws.Range("E12").Value = "  +1.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.848.82"
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  +2.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.487"
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.24"
$ws.Range("D16").Style = $defaultStyle
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008702"
$ws.Range("D18").Style = $defaultStyle
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.014"
$ws.Range("D19").Style = $defaultStyle
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.143.94"
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.58"
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = "  +2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.057"
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.944"
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.00"
$ws.Range("D25").Style = $defaultStyle
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.36"
$ws.Range("D26").Style = $defaultStyle
$ws.Range("E26").Value = "  +1.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.986"
$ws.Range("D27").Style = $defaultStyle
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.74"
$ws.Range("D28").Style = $defaultStyle
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.921"
$ws.Range("D29").Style = $defaultStyle
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08864"
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.317"
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = "  +3.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.182"
$ws.Range("D32").Style = $defaultStyle
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.522"
$ws.Range("D33").Style = $defaultStyle
$ws.Range("E33").Value = "  +1.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7389"
$ws.Range("D34").Style = $defaultStyle
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.690"
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = "  -2.97%  "
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01975"
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05266"
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = "  +2.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.973"
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5280"
$ws.Range("D40").Style = $defaultStyle
$ws.Range("E40").Value = "  +2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.043"
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = "  +2.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1524"
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.296"
$ws.Range("D43").Style = $defaultStyle
$ws.Range("E43").Value = "  +2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.61"
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = "  +1.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4745"
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.07"
$ws.Range("D47").Style = $defaultStyle
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.616"
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = "  +1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "66.07"
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06059"
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8931"
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = "  +3.53%  "
